$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("API")
$ws.Range("B2").Value = 0.4061
$ws.Range("C2").Value = 31686.6893
$ws.Range("D2").Value = 1
$ws.Range("B3").Value = -0.2594
$ws.Range("C3").Value = 0.2233
$ws.Range("D3").Value = 0.2453
$ws.Range("B4").Value = -0.1332
$ws.Range("C4").Value = 0.1692
$ws.Range("D4").Value = 0.4313
$ws.Range("B5").Value = 1.6647
$ws.Range("C5").Value = 31686.6893
$ws.Range("D5").Value = 1
$ws.Range("B6").Value = -0.2284
$ws.Range("C6").Value = 40824.4892
$ws.Range("D6").Value = 1
$ws.Range("B7").Value = -0.2215
$ws.Range("C7").Value = 38695.3832
$ws.Range("D7").Value = 1
$ws.Range("B8").Value = -0.1344
$ws.Range("C8").Value = 42715.2131
$ws.Range("D8").Value = 1
$ws.Range("B9").Value = 1.19
$ws.Range("C9").Value = 31686.6893
$ws.Range("D9").Value = 1
$ws.Range("B10").Value = 1.5852
$ws.Range("C10").Value = 31686.6893
$ws.Range("D10").Value = 1
$ws.Range("B11").Value = -0.0641
$ws.Range("C11").Value = 42580.1311
$ws.Range("D11").Value = 1
$ws.Range("B12").Value = -0.5778
$ws.Range("C12").Value = 0.6369
$ws.Range("D12").Value = 0.3643
$ws.Range("B13").Value = -1.0665
$ws.Range("C13").Value = 0.8657
$ws.Range("D13").Value = 0.218
$ws.Range("B14").Value = -1.2398
$ws.Range("C14").Value = 0.8763
$ws.Range("D14").Value = 0.1571
$ws.Range("B15").Value = -1.1869
$ws.Range("C15").Value = 0.8084
$ws.Range("D15").Value = 0.1421
$ws.Range("B16").Value = -1.0491
$ws.Range("C16").Value = 0.7556
$ws.Range("D16").Value = 0.165
$ws.Range("B17").Value = -0.8321
$ws.Range("C17").Value = 0.7133
$ws.Range("D17").Value = 0.2434
$ws.Range("B18").Value = -0.7644
$ws.Range("C18").Value = 0.6983
$ws.Range("D18").Value = 0.2736
$ws.Range("B19").Value = -0.9932
$ws.Range("C19").Value = 0.6979
$ws.Range("D19").Value = 0.1547
$ws.Range("B20").Value = -1.1157
$ws.Range("C20").Value = 0.7571
$ws.Range("D20").Value = 0.1406
$ws.Range("B21").Value = -1.1348
$ws.Range("C21").Value = 0.8509
$ws.Range("D21").Value = 0.1823
$ws.Range("B22").Value = -1.1315
$ws.Range("C22").Value = 1.151
$ws.Range("D22").Value = 0.3256
$ws.Range("B23").Value = -1.1283
$ws.Range("C23").Value = 1.7581
$ws.Range("D23").Value = 0.521
$ws.Range("B24").Value = -1.225
$ws.Range("C24").Value = 0.2217

$ws = $wb.Worksheets.Item("POP")
$ws.Range("B2").Value = 2.315
$ws.Range("C2").Value = 0.2838
$ws.Range("B3").Value = -0.2653
$ws.Range("C3").Value = 0.0817
$ws.Range("D3").Value = 0.0012
$ws.Range("B4").Value = -0.3223
$ws.Range("C4").Value = 0.0884
$ws.Range("D4").Value = 0.0003
$ws.Range("B5").Value = -0.3065
$ws.Range("C5").Value = 0.2529
$ws.Range("D5").Value = 0.2254
$ws.Range("B6").Value = 0.2063
$ws.Range("C6").Value = 0.1625
$ws.Range("D6").Value = 0.2041
$ws.Range("B7").Value = 0.061
$ws.Range("C7").Value = 0.1622
$ws.Range("D7").Value = 0.707
$ws.Range("B8").Value = -0.4051
$ws.Range("C8").Value = 0.342
$ws.Range("D8").Value = 0.2362
$ws.Range("B9").Value = 0.3051
$ws.Range("C9").Value = 0.1764
$ws.Range("D9").Value = 0.0836
$ws.Range("B10").Value = 0.2152
$ws.Range("C10").Value = 0.1585
$ws.Range("D10").Value = 0.1745
$ws.Range("B11").Value = -0.4734
$ws.Range("C11").Value = 0.5145
$ws.Range("D11").Value = 0.3575
$ws.Range("B12").Value = -0.0528
$ws.Range("C12").Value = 0.1969
$ws.Range("D12").Value = 0.7886
$ws.Range("B13").Value = -0.0986
$ws.Range("C13").Value = 0.2659
$ws.Range("D13").Value = 0.7109
$ws.Range("B14").Value = -0.1032
$ws.Range("C14").Value = 0.2718
$ws.Range("D14").Value = 0.7043
$ws.Range("B15").Value = -0.0517
$ws.Range("C15").Value = 0.2678
$ws.Range("D15").Value = 0.847
$ws.Range("B16").Value = 0.0129
$ws.Range("C16").Value = 0.2559
$ws.Range("D16").Value = 0.9599
$ws.Range("B17").Value = 0.0513
$ws.Range("C17").Value = 0.2506
$ws.Range("D17").Value = 0.8379
$ws.Range("B18").Value = 0.0151
$ws.Range("C18").Value = 0.2487
$ws.Range("D18").Value = 0.9516
$ws.Range("B19").Value = -0.0803
$ws.Range("C19").Value = 0.2484
$ws.Range("D19").Value = 0.7466
$ws.Range("B20").Value = -0.1165
$ws.Range("C20").Value = 0.2637
$ws.Range("D20").Value = 0.6586
$ws.Range("B21").Value = -0.0855
$ws.Range("C21").Value = 0.2735
$ws.Range("D21").Value = 0.7547
$ws.Range("B22").Value = -0.0216
$ws.Range("C22").Value = 0.2851
$ws.Range("D22").Value = 0.9396
$ws.Range("B23").Value = 0.048
$ws.Range("C23").Value = 0.3831
$ws.Range("D23").Value = 0.9003
$ws.Range("B24").Value = -1.3513
$ws.Range("C24").Value = 0.0814

$ws = $wb.Worksheets.Item("Pesticide")
$ws.Range("B2").Value = 3.3025
$ws.Range("C2").Value = 1.2803
$ws.Range("D2").Value = 0.0099
$ws.Range("B3").Value = 0.5262
$ws.Range("C3").Value = 0.243
$ws.Range("D3").Value = 0.0303
$ws.Range("B4").Value = 0.3542
$ws.Range("C4").Value = 0.233
$ws.Range("D4").Value = 0.1285
$ws.Range("B5").Value = 0.607
$ws.Range("C5").Value = 0.5845
$ws.Range("D5").Value = 0.299
$ws.Range("B6").Value = -0.2895
$ws.Range("C6").Value = 0.7234
$ws.Range("D6").Value = 0.6891
$ws.Range("B7").Value = -2.0377
$ws.Range("C7").Value = 2189.6613
$ws.Range("D7").Value = 0.9993
$ws.Range("B8").Value = 0.6646
$ws.Range("C8").Value = 0.5659
$ws.Range("D8").Value = 0.2402
$ws.Range("B9").Value = -0.4127
$ws.Range("C9").Value = 0.7752
$ws.Range("D9").Value = 0.5945
$ws.Range("B10").Value = 0.1365
$ws.Range("C10").Value = 0.6849
$ws.Range("D10").Value = 0.842
$ws.Range("B11").Value = -1.4508
$ws.Range("C11").Value = 2981.6323
$ws.Range("D11").Value = 0.9996
$ws.Range("B12").Value = -1.2623
$ws.Range("C12").Value = 0.9743
$ws.Range("D12").Value = 0.1951
$ws.Range("B13").Value = -2.3668
$ws.Range("C13").Value = 1.4649
$ws.Range("D13").Value = 0.1062
$ws.Range("B14").Value = -2.9957
$ws.Range("C14").Value = 1.7185
$ws.Range("D14").Value = 0.0813
$ws.Range("B15").Value = -3.2126
$ws.Range("C15").Value = 1.744
$ws.Range("D15").Value = 0.0655
$ws.Range("B16").Value = -3.1079
$ws.Range("C16").Value = 1.5638
$ws.Range("D16").Value = 0.0469
$ws.Range("B17").Value = -2.7713
$ws.Range("C17").Value = 1.3013
$ws.Range("D17").Value = 0.0332
$ws.Range("B18").Value = -2.2698
$ws.Range("C18").Value = 1.1417
$ws.Range("D18").Value = 0.0468
$ws.Range("B19").Value = -2.1592
$ws.Range("C19").Value = 1.1317
$ws.Range("D19").Value = 0.0564
$ws.Range("B20").Value = -2.2244
$ws.Range("C20").Value = 1.1697
$ws.Range("D20").Value = 0.0572
$ws.Range("B21").Value = -2.0845
$ws.Range("C21").Value = 1.187
$ws.Range("D21").Value = 0.0791
$ws.Range("B22").Value = -2.3118
$ws.Range("C22").Value = 1.2185
$ws.Range("D22").Value = 0.0578
$ws.Range("B23").Value = -2.6108
$ws.Range("C23").Value = 1.5281
$ws.Range("D23").Value = 0.0875
$ws.Range("B24").Value = -0.9765
$ws.Range("C24").Value = 0.2063
